$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Well settings")
